$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (regenerated s_val data, filtering save games)
$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 3099.503889238888
$ws.Range("D2").Value = 337.1190423067083
$ws.Range("E2").Value = 616238.5361209477
$ws.Range("G2").Value = 619675.8343826485

# Row 3 values
$ws.Range("B3").Value = 0.04763786555579896
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 616238.5361209477
$ws.Range("G3").Value = 616265.1139556493
